$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper functions -------------------------------------------------
# Set a cell so it becomes a shared "text" cell (e.g. "0" or "***.*"),
# using the style of a known-good reference cell that already has that
# display ("right, General format" text style).
function Set-TextCell {
    param($cellRef, $text, $styleSourceRef)
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($styleSourceRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# Set a cell so it becomes a numeric cell with the style of a known-good
# reference cell (number or percent format).
function Set-NumberCell {
    param($cellRef, $value, $styleSourceRef)
    $ws.Range($styleSourceRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($cellRef).Value2 = $value
}

# --- Shared-string text updates (report header) ------------------------
# "Volume 31   Number  2" -> "Volume 31   Number  3"
$ws.Range("A8").Characters(21, 1).Text = "3"

# "Report Covering the Week  1/8/2024  Through  1/14/2024"
#  -> "Report Covering the Week  1/15/2024  Through  1/21/2024"
$ws.Range("C9").Characters(27, 8).Text = "1/15/2024"
$ws.Range("C9").Characters(46, 9).Text = "1/21/2024"

# --- Row 15 (Rape) ------------------------------------------------------
Set-TextCell "D15" "0" "C23"
Set-TextCell "E15" "***.*" "E23"
Set-NumberCell "N15" -100 "K36"

# --- Row 16 (Robbery) ----------------------------------------------------
$ws.Range("C16").Value2 = 4
$ws.Range("D16").Value2 = 3
$ws.Range("E16").Value2 = 33.333333333333
$ws.Range("F16").Value2 = 13
$ws.Range("G16").Value2 = 7
$ws.Range("H16").Value2 = 85.714285714285
$ws.Range("I16").Value2 = 8
$ws.Range("J16").Value2 = 7
$ws.Range("K16").Value2 = 14.285714285714
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -27.272727272727
$ws.Range("N16").Value2 = -86.206896551724

# --- Row 17 (Fel. Assault) -------------------------------------------------
$ws.Range("C17").Value2 = 6
$ws.Range("D17").Value2 = 10
$ws.Range("E17").Value2 = -40
$ws.Range("G17").Value2 = 19
$ws.Range("H17").Value2 = 5.263157894736
$ws.Range("I17").Value2 = 17
$ws.Range("J17").Value2 = 15
$ws.Range("K17").Value2 = 13.333333333333
$ws.Range("L17").Value2 = 183.333333333333
$ws.Range("M17").Value2 = 466.666666666667
$ws.Range("N17").Value2 = 30.769230769230

# --- Row 18 (Burglary) ---------------------------------------------------
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 3
$ws.Range("E18").Value2 = 0
$ws.Range("I18").Value2 = 5
$ws.Range("J18").Value2 = 7
$ws.Range("K18").Value2 = -28.571428571428
$ws.Range("L18").Value2 = 25
$ws.Range("M18").Value2 = -80.769230769230
$ws.Range("N18").Value2 = -94.791666666666

# --- Row 19 (Gr. Larceny) -------------------------------------------------
$ws.Range("C19").Value2 = 12
$ws.Range("D19").Value2 = 9
$ws.Range("E19").Value2 = 33.333333333333
$ws.Range("F19").Value2 = 51
$ws.Range("G19").Value2 = 38
$ws.Range("H19").Value2 = 34.210526315789
$ws.Range("I19").Value2 = 37
$ws.Range("J19").Value2 = 31
$ws.Range("K19").Value2 = 19.354838709677
$ws.Range("L19").Value2 = -40.322580645161
$ws.Range("M19").Value2 = 32.142857142857
$ws.Range("N19").Value2 = -9.756097560975

# --- Row 20 (G.L.A.) -------------------------------------------------------
$ws.Range("C20").Value2 = 4
$ws.Range("E20").Value2 = 100
$ws.Range("F20").Value2 = 13
$ws.Range("G20").Value2 = 8
$ws.Range("H20").Value2 = 62.5
$ws.Range("I20").Value2 = 11
$ws.Range("J20").Value2 = 4
$ws.Range("K20").Value2 = 175
$ws.Range("L20").Value2 = 83.333333333333
$ws.Range("M20").Value2 = 10
$ws.Range("N20").Value2 = -89.622641509434

# --- Row 21 (TOTAL) --------------------------------------------------------
$ws.Range("C21").Value2 = 29
$ws.Range("D21").Value2 = 27
$ws.Range("E21").Value2 = 7.407407407407
$ws.Range("F21").Value2 = 108
$ws.Range("G21").Value2 = 84
$ws.Range("H21").Value2 = 28.571428571428
$ws.Range("I21").Value2 = 78
$ws.Range("J21").Value2 = 67
$ws.Range("K21").Value2 = 16.417910447761
$ws.Range("L21").Value2 = -11.363636363636
$ws.Range("M21").Value2 = -1.265822784810
$ws.Range("N21").Value2 = -75.316455696202

# --- Row 22 (Transit) -------------------------------------------------------
Set-NumberCell "C22" 1 "C36"
$ws.Range("F22").Value2 = 2
$ws.Range("I22").Value2 = 2
Set-NumberCell "L22" 100 "K36"

# --- Row 24 (Petit Larceny) --------------------------------------------------
$ws.Range("C24").Value2 = 39
$ws.Range("D24").Value2 = 36
$ws.Range("E24").Value2 = 8.333333333333
$ws.Range("F24").Value2 = 139
$ws.Range("G24").Value2 = 113
$ws.Range("H24").Value2 = 23.008849557522
$ws.Range("I24").Value2 = 111
$ws.Range("J24").Value2 = 84
$ws.Range("K24").Value2 = 32.142857142857
$ws.Range("L24").Value2 = 13.265306122449
$ws.Range("M24").Value2 = 88.135593220339

# --- Row 25 (Misd. Assault) ---------------------------------------------------
$ws.Range("C25").Value2 = 14
$ws.Range("D25").Value2 = 10
$ws.Range("E25").Value2 = 40
$ws.Range("F25").Value2 = 65
$ws.Range("G25").Value2 = 37
$ws.Range("H25").Value2 = 75.675675675675
$ws.Range("I25").Value2 = 48
$ws.Range("J25").Value2 = 26
$ws.Range("K25").Value2 = 84.615384615384
$ws.Range("L25").Value2 = 108.695652173913
$ws.Range("M25").Value2 = 118.181818181818

# --- Row 26 (UCR Rape*) --------------------------------------------------------
Set-TextCell "D26" "0" "C23"
Set-TextCell "E26" "***.*" "E23"
Set-NumberCell "F26" 1 "C36"
$ws.Range("H26").Value2 = -75
Set-NumberCell "I26" 1 "C36"
$ws.Range("K26").Value2 = -66.666666666666
$ws.Range("L26").Value2 = -50

# --- Row 27 (Other Sex Crimes) --------------------------------------------------
$ws.Range("D27").Value2 = 1
$ws.Range("E27").Value2 = 0
$ws.Range("J27").Value2 = 3
$ws.Range("K27").Value2 = 66.666666666666
$ws.Range("L27").Value2 = 150

# --- Row 30 (Hate Crimes) ---------------------------------------------------------
Set-TextCell "F30" "0" "C23"
